$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (pushes existing rows 14-83 down to 15-84)
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with its data.
$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(14, 3).Value = "Bíobío"
$ws.Cells.Item(14, 4).Value = 44670
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 100112001
$ws.Cells.Item(14, 7).Value = "Berenjena"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 180
$ws.Cells.Item(14, 11).Value = 6000
$ws.Cells.Item(14, 12).Value = 7000
$ws.Cells.Item(14, 13).Value = 6444
$ws.Cells.Item(14, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(14, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(14, 16).Value = 107
$ws.Cells.Item(14, 17).Value = 60
$ws.Cells.Item(14, 18).Value = "Hortaliza"
